$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Range("L2").Value = 6269
$ws.Range("L3").Value = 6795
$ws.Range("I4").Value = 1851
$ws.Range("L4").Value = 1688
$ws.Range("L5").Value = 402
$ws.Range("L6").Value = 5585
$ws.Range("I7").Value = 26321
$ws.Range("L7").Value = 20739

$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Range("L2").Value = 184
$ws.Range("L4").Value = 75
$ws.Range("L7").Value = 659
$ws.Range("L8").Value = 1373
$ws.Range("L11").Value = 342
$ws.Range("L18").Value = 143
$ws.Range("L19").Value = 561
$ws.Range("L20").Value = 527
$ws.Range("L22").Value = 69
$ws.Range("L25").Value = 126
$ws.Range("L29").Value = 1156
$ws.Range("L31").Value = 206
$ws.Range("L33").Value = 935
$ws.Range("L37").Value = 796
$ws.Range("L42").Value = 659
$ws.Range("L43").Value = 156
$ws.Range("L48").Value = 273
$ws.Range("L51").Value = 258
$ws.Range("L52").Value = 442
$ws.Range("L54").Value = 445
$ws.Range("I63").Value = 271
$ws.Range("L65").Value = 405
$ws.Range("L67").Value = 722
$ws.Range("L70").Value = 37
$ws.Range("L76").Value = 322
$ws.Range("L77").Value = 138
$ws.Range("L78").Value = 274
$ws.Range("L79").Value = 571
$ws.Range("L80").Value = 69
$ws.Range("L83").Value = 454
$ws.Range("L84").Value = 200
$ws.Range("L85").Value = 1029
$ws.Range("L90").Value = 219
$ws.Range("L96").Value = 230
$ws.Range("L97").Value = 167
$ws.Range("I101").Value = 26321
$ws.Range("L101").Value = 20739

$ws = $wb.Worksheets.Item('West Ridge')
$ws.Range("L4").Value = 24
$ws.Range("L7").Value = 230

$ws = $wb.Worksheets.Item('Auburn Gresham')
$ws.Range("L4").Value = 45
$ws.Range("L7").Value = 659

$ws = $wb.Worksheets.Item('Belmont Cragin')
$ws.Range("L5").Value = 3
$ws.Range("L7").Value = 342

$ws = $wb.Worksheets.Item('South Shore')
$ws.Range("L2").Value = 310
$ws.Range("L4").Value = 58
$ws.Range("L6").Value = 212
$ws.Range("L7").Value = 1029

$ws = $wb.Worksheets.Item('Little Village')
$ws.Range("L6").Value = 126
$ws.Range("L7").Value = 442

$ws = $wb.Worksheets.Item('Austin')
$ws.Range("L2").Value = 418
$ws.Range("L6").Value = 332
$ws.Range("L7").Value = 1373

$ws = $wb.Worksheets.Item('South Chicago')
$ws.Range("L6").Value = 102
$ws.Range("L7").Value = 454

$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Range("L2").Value = 254
$ws.Range("L3").Value = 330
$ws.Range("L7").Value = 935

$ws = $wb.Worksheets.Item('Grand Crossing')
$ws.Range("L3").Value = 283
$ws.Range("L4").Value = 44
$ws.Range("L7").Value = 796

$ws = $wb.Worksheets.Item('New City')
$ws.Range("L2").Value = 147
$ws.Range("L7").Value = 405

$ws = $wb.Worksheets.Item('Gage Park')
$ws.Range("L2").Value = 81
$ws.Range("L4").Value = 12
$ws.Range("L7").Value = 206

$ws = $wb.Worksheets.Item('North Lawndale')
$ws.Range("L2").Value = 205
$ws.Range("L3").Value = 281
$ws.Range("L7").Value = 722

$ws = $wb.Worksheets.Item('South Deering')
$ws.Range("L2").Value = 68
$ws.Range("L7").Value = 200

$ws = $wb.Worksheets.Item('Loop')
$ws.Range("L4").Value = 37
$ws.Range("L7").Value = 445

$ws = $wb.Worksheets.Item('Englewood')
$ws.Range("L3").Value = 445
$ws.Range("L6").Value = 282
$ws.Range("L7").Value = 1156

$ws = $wb.Worksheets.Item('Lake View')
$ws.Range("L3").Value = 70
$ws.Range("L7").Value = 273

$ws = $wb.Worksheets.Item('Chatham')
$ws.Range("L6").Value = 153
$ws.Range("L7").Value = 561

$ws = $wb.Worksheets.Item('River North')
$ws.Range("L2").Value = 66
$ws.Range("L4").Value = 39
$ws.Range("L7").Value = 322

$ws = $wb.Worksheets.Item('Humboldt Park')
$ws.Range("L6").Value = 187
$ws.Range("L7").Value = 659

$ws = $wb.Worksheets.Item('Rogers Park')
$ws.Range("L3").Value = 91
$ws.Range("L7").Value = 274

$ws = $wb.Worksheets.Item('Roseland')
$ws.Range("L4").Value = 41
$ws.Range("L7").Value = 571

$ws = $wb.Worksheets.Item('Chicago Lawn')
$ws.Range("L2").Value = 165
$ws.Range("L7").Value = 527

$ws = $wb.Worksheets.Item('Calumet Heights')
$ws.Range("L2").Value = 50
$ws.Range("L5").Value = 4
$ws.Range("L7").Value = 143

$ws = $wb.Worksheets.Item('East Side')
$ws.Range("L6").Value = 19
$ws.Range("L7").Value = 126

$ws = $wb.Worksheets.Item('Albany Park')
$ws.Range("L3").Value = 60
$ws.Range("L7").Value = 184

$ws = $wb.Worksheets.Item('West Town')
$ws.Range("L3").Value = 37
$ws.Range("L7").Value = 167

$ws = $wb.Worksheets.Item('O''Hare')
$ws.Range("L3").Value = 18
$ws.Range("L7").Value = 37

$ws = $wb.Worksheets.Item('Edgewater')
$ws.Range("L3").Value = 50
$ws.Range("L6").Value = 54

$ws = $wb.Worksheets.Item('Washington Heights')
$ws.Range("L3").Value = 63
$ws.Range("L7").Value = 219

$ws = $wb.Worksheets.Item('Little Italy, UIC')
$ws.Range("L6").Value = 55
$ws.Range("L7").Value = 258

$ws = $wb.Worksheets.Item('Hyde Park')
$ws.Range("L2").Value = 32
$ws.Range("L4").Value = 26
$ws.Range("L7").Value = 156

$ws = $wb.Worksheets.Item('Clearing')
$ws.Range("L3").Value = 28
$ws.Range("L7").Value = 69

$ws = $wb.Worksheets.Item('Riverdale')
$ws.Range("L3").Value = 43
$ws.Range("L7").Value = 138

$ws = $wb.Worksheets.Item('Rush & Division')
$ws.Range("L6").Value = 32
$ws.Range("L7").Value = 69

$ws = $wb.Worksheets.Item('Archer Heights')
$ws.Range("L3").Value = 19
$ws.Range("L7").Value = 75
